# Automatische test-sync: 2025-06-29 15:06:50
# Appends the 12th test-mail log entry ("Offerte / Prijsaanvraag") to the
# "Logs" sheet, rolls the corresponding tally up on the "Dashboard" sheet,
# and extends the conditional-formatting ranges + bar-chart series so they
# keep covering the newly-added row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 27 with the new test-mail entry
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(27, 1).Value = "Wanneer komt mijn offerte?"
$logs.Cells.Item(27, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(27, 3).Value = "Testmail #12: Wanneer komt mijn offerte?"
$logs.Cells.Item(27, 4).Value = "Offerte / Prijsaanvraag"
$logs.Cells.Item(27, 5).Value = "Geachte klant,`nDank u voor uw e-mail. Uw offerte zal naar verwachting binnen 24 uur worden verstuurd. Mocht u deze niet op tijd ontvangen, neem dan gerust contact met ons op.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Cells.Item(27, 6).Value = "2025-06-29 15:06:42"
$logs.Cells.Item(27, 7).Value = "Ja"
$logs.Cells.Item(27, 8).Value = "Nee"
$logs.Cells.Item(27, 9).Value = "Ja"

# Extend the conditional-formatting ranges (D/G/H/I 2:26 -> 2:27) so the
# newly-added row is covered, same as Excel does when you drag/copy the
# formatting down with the rest of the table.
$cfColumns = @("D", "G", "H", "I")
foreach ($col in $cfColumns) {
    $oldRange = $logs.Range($col + "2:" + $col + "26")
    $newRange = $logs.Range($col + "2:" + $col + "27")
    $conditions = $oldRange.FormatConditions
    for ($i = 1; $i -le $conditions.Count; $i++) {
        $conditions.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 2. Dashboard sheet: append the new category tally row 9
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Cells.Item(9, 1).Value = "Offerte / Prijsaanvraag"
$dashboard.Cells.Item(9, 2).Value = 1

# ---------------------------------------------------------------------
# 3. Chart: extend the category/value series ranges from row 8 to row 9
# ---------------------------------------------------------------------
$chartObj = $dashboard.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$9,Dashboard!`$B`$2:`$B`$9,1)"
